$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 ("Recorded By")
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
